$wb = $excel.ActiveWorkbook

# --- Update control-word bits on Foglio3 for the jump instructions ---
# Columns I, J, K, L correspond to instructions: beqz, bnez, j, jal
$ws3 = $wb.Worksheets.Item("Foglio3")

# Row 14 = sel_op1_mux: was 1 for these jump instrs, now 0
$ws3.Range("I14:L14").Value = 0

# Row 16 = ALU_func_0, Row 17 = ALU_func_1, Row 18 = ALU_func_2, Row 19 = ALU_func_3
# were 0 for these jump instrs, now 1
$ws3.Range("I16:L16").Value = 1
$ws3.Range("I17:L17").Value = 1
$ws3.Range("I18:L18").Value = 1
$ws3.Range("I19:L19").Value = 1

# --- Update selections / active sheet to match the saved view state ---
$ws3.Activate()
$ws3.Range("K22").Select()

$ws2 = $wb.Worksheets.Item("Foglio2")
$ws2.Activate()
$ws2.Range("E9").Select()
